# Updates cryptos list data (prices and 1h volume % changes) as described
# by the commit message: "Updated cryptos list on Tue Nov 14 04:55:46 UTC
# 2023 with GitHub Actions". Rows 30/31 also swap coins (ImmutableX /
# Stellar) to reflect the refreshed ranking order/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-ambiguous values: assign directly. ---
$ws.Range('D2').Value = '36.524.90'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '2.056.13'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('E5').Value = '  -0.82%  '
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -5.44%  '
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('E10').Value = '  -3.89%  '
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('E12').Value = '  -2.83%  '
$ws.Range('E13').Value = '  +4.64%  '
$ws.Range('E14').Value = '  -4.26%  '
$ws.Range('D15').Value = '2.357.62'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('E16').Value = '  -3.93%  '
$ws.Range('D17').Value = '2.051.20'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '36.493.52'
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('E19').Value = '  -6.99%  '
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('E21').Value = '  -3.14%  '
$ws.Range('E22').Value = '  +1.40%  '
$ws.Range('E23').Value = '  -2.19%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('E25').Value = '  -3.75%  '
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('E28').Value = '  -2.47%  '
$ws.Range('E29').Value = '  +1.24%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E31').Value = '  +11.79%  '
$ws.Range('E32').Value = '  -4.98%  '
$ws.Range('E33').Value = '  -4.03%  '
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('E38').Value = '  -5.87%  '
$ws.Range('E39').Value = '  -4.41%  '
$ws.Range('E40').Value = '  -4.97%  '
$ws.Range('E41').Value = '  -2.48%  '
$ws.Range('E42').Value = '  -7.86%  '
$ws.Range('E43').Value = '  -5.10%  '
$ws.Range('E44').Value = '  -2.31%  '
$ws.Range('E45').Value = '  -2.65%  '
$ws.Range('D46').Value = '1.416.81'
$ws.Range('E46').Value = '  +9.64%  '
$ws.Range('E47').Value = '  +13.43%  '
$ws.Range('E48').Value = '  -5.11%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  -2.39%  '
$ws.Range('D51').Value = '2.247.95'
$ws.Range('E51').Value = '  +0.91%  '

# --- Price values that look like plain numbers (e.g. '242.86') need to
# stay text-typed (matching the source sheet's inlineStr cells) instead of
# being auto-coerced to a numeric cell by a direct .Value assignment. Route
# them through a quoted formula then Copy/PasteSpecial values-only, which
# collapses the formula down to a literal text cell with no style churn. ---
$priceCells = [ordered]@{
    'D5' = '242.86'
    'D8' = '54.64'
    'D9' = '58.43'
    'D13' = '0.916'
    'D14' = '14.72'
    'D16' = '5.40'
    'D19' = '16.81'
    'D20' = '71.99'
    'D22' = '238.67'
    'D23' = '5.24'
    'D25' = '2.35'
    'D26' = '9.33'
    'D28' = '164.96'
    'D29' = '20.07'
    'D30' = '0.122'
    'D31' = '1.22'
    'D32' = '5.09'
    'D33' = '4.47'
    'D36' = '1.84'
    'D38' = '0.0821'
    'D39' = '1.24'
    'D42' = '2.86'
    'D43' = '0.0926'
    'D45' = '93.83'
    'D47' = '7.60'
    'D48' = '15.91'
}
foreach ($ref in $priceCells.Keys) {
    $ws.Range($ref).Formula = '="' + $priceCells[$ref] + '"'
    $ws.Range($ref).Copy()
    $ws.Range($ref).PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
